$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H70").Value = 6315.4546
$ws.Range("I70").Value = 11618.8
$ws.Range("J70").Value = 1896
$ws.Range("K70").Value = 34856.39999999999
$ws.Range("L70").Value = 5688
$ws.Range("M70").Value = -34586.39999999999
$ws.Range("N70").Value = -6228
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H73").Value = 6315.4546
$ws.Range("I73").Value = 11618.8
$ws.Range("J73").Value = 1896
$ws.Range("K73").Value = 34856.39999999999
$ws.Range("L73").Value = 5688
$ws.Range("M73").Value = -33920.39999999999
$ws.Range("N73").Value = -7560
$ws.Range("H86").Value = 2147.9092
$ws.Range("J86").Value = 1999.8
$ws.Range("L86").Value = 1999.8
$ws.Range("N86").Value = -4245.8
$ws.Range("H89").Value = 2147.9092
$ws.Range("J89").Value = 1999.8
$ws.Range("L89").Value = 9999
$ws.Range("N89").Value = -21231
$ws.Range("H103").Value = 415.17856
$ws.Range("J103").Value = 451.13635
$ws.Range("L103").Value = 1353.40905
$ws.Range("N103").Value = -2525.40905
$ws.Range("H106").Value = 11818.091
$ws.Range("I106").Value = 2167.5
$ws.Range("K106").Value = 2167.5
$ws.Range("M106").Value = -1536.5
$ws.Range("H135").Value = 34072.066
$ws.Range("I135").Value = 694.8889
$ws.Range("K135").Value = 6254.0001
$ws.Range("M135").Value = -3719.0001
$ws.Range("H137").Value = 33626.13
$ws.Range("J137").Value = 72415.78999999999
$ws.Range("L137").Value = 217247.37
$ws.Range("N137").Value = -222347.37

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 19900
$ws.Range("H80").Value = 19998
$ws.Range("J80").Value = 19998
$ws.Range("L80").Value = 19998
$ws.Range("N80").Value = -21994
$ws.Range("H83").Value = 19998
$ws.Range("J83").Value = 19998
$ws.Range("L83").Value = 59994
$ws.Range("N83").Value = -69978
$ws.Range("H104").Value = 30224
$ws.Range("J104").Value = 30224
$ws.Range("L104").Value = 30224
$ws.Range("N104").Value = -37212

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1135.8422
$ws.Range("I80").Value = 847.2222
$ws.Range("J80").Value = 1395.6
$ws.Range("K80").Value = 847.2222
$ws.Range("L80").Value = 1395.6
$ws.Range("M80").Value = 150.7778
$ws.Range("N80").Value = -3391.6
$ws.Range("H83").Value = 1135.8422
$ws.Range("I83").Value = 847.2222
$ws.Range("J83").Value = 1395.6
$ws.Range("K83").Value = 4236.111
$ws.Range("L83").Value = 6978
$ws.Range("M83").Value = 755.8890000000001
$ws.Range("N83").Value = -16962
$ws.Range("H94").Value = 1047.25
$ws.Range("I94").Value = 1047.25
$ws.Range("K94").Value = 1047.25
$ws.Range("M94").Value = -596.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1079
$ws.Range("J16").Value = 1300
$ws.Range("L16").Value = 1300
$ws.Range("N16").Value = -1874
$ws.Range("H62").Value = 27377.75
$ws.Range("I62").Value = 35835
$ws.Range("K62").Value = 35835
$ws.Range("M62").Value = -35211
$ws.Range("H65").Value = 27377.75
$ws.Range("I65").Value = 35835
$ws.Range("K65").Value = 179175
$ws.Range("M65").Value = -176055
$ws.Range("H74").Value = 33933.332
$ws.Range("J74").Value = 33933.332
$ws.Range("L74").Value = 33933.332
$ws.Range("N74").Value = -35681.332
$ws.Range("H77").Value = 33933.332
$ws.Range("J77").Value = 33933.332
$ws.Range("L77").Value = 101799.996
$ws.Range("N77").Value = -110535.996
$ws.Range("H94").Value = 2012.6666
$ws.Range("I94").Value = 1321
$ws.Range("K94").Value = 1321
$ws.Range("M94").Value = -870
$ws.Range("H105").Value = 6582
$ws.Range("I105").Value = 6582
$ws.Range("K105").Value = 6582
$ws.Range("M105").Value = -4835
$ws.Range("H113").Value = 1079
$ws.Range("J113").Value = 1300
$ws.Range("L113").Value = 1300
$ws.Range("N113").Value = -5640
$ws.Range("H122").Value = 2671.4092
$ws.Range("J122").Value = 1817.5555
$ws.Range("L122").Value = 5452.666499999999
$ws.Range("N122").Value = -10352.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 15174.5
$ws.Range("I99").Value = 6349
$ws.Range("K99").Value = 19047
$ws.Range("M99").Value = -16801

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 5393.857
$ws.Range("I19").Value = 6771
$ws.Range("J19").Value = 4361
$ws.Range("K19").Value = 6771
$ws.Range("L19").Value = 4361
$ws.Range("M19").Value = -6483
$ws.Range("N19").Value = -4937
$ws.Range("H104").Value = 29171
$ws.Range("J104").Value = 29171
$ws.Range("L104").Value = 29171
$ws.Range("N104").Value = -36159
$ws.Range("H123").Value = 63993.5
$ws.Range("J123").Value = 63993.5
$ws.Range("L123").Value = 63993.5
$ws.Range("N123").Value = -68893.5
$ws.Range("H136").Value = 40382.91
$ws.Range("J136").Value = 40382.91
$ws.Range("L136").Value = 121148.73
$ws.Range("N136").Value = -126248.73

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1199.75
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1199.75
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1199.75
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1789.75
$ws.Range("H27").Value = 1199.75
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1199.75
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1199.75
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1413.75
$ws.Range("H133").Value = 82665
$ws.Range("J133").Value = 82665
$ws.Range("L133").Value = 82665
$ws.Range("N133").Value = -87725

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H81").Value = 4800.375
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 4800.375
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H105").Value = 26410
$ws.Range("J105").Value = 26410
$ws.Range("L105").Value = 26410
$ws.Range("N105").Value = -33398
$ws.Range("H133").Value = 46854.43
$ws.Range("J133").Value = 46854.43
$ws.Range("L133").Value = 46854.43
$ws.Range("N133").Value = -56974.43
$ws.Range("H136").Value = 699.82355
$ws.Range("I136").Value = 706.30304
$ws.Range("J136").Value = 486
$ws.Range("K136").Value = 2118.90912
$ws.Range("L136").Value = 1458
$ws.Range("M136").Value = 431.0908799999997
$ws.Range("N136").Value = -6558
